$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value2 = "SAMPLE_TYPE"
$ws.Range("A1").Value2 = "'Specimen_Number"

$ws.Range("A1:G1").Select() | Out-Null
